$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 38, pushing existing rows 38-44 down to 40-46.
$ws.Rows("38:39").Insert()

# Row 38 : new "Carson" record
$ws.Cells.Item(38,1).Value = 1
$ws.Cells.Item(38,2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(38,3).Value = "Arica y Parinacota"
$ws.Cells.Item(38,4).Value = 44615
$ws.Cells.Item(38,5).Value = 15
$ws.Cells.Item(38,6).Value = "Fruta"
$ws.Cells.Item(38,7).Value = 100103
$ws.Cells.Item(38,8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(38,9).Value = 100103004
$ws.Cells.Item(38,10).Value = "Durazno"
$ws.Cells.Item(38,11).Value = "Carson"
$ws.Cells.Item(38,12).Value = "Primera"
$ws.Cells.Item(38,13).Value = 270
$ws.Cells.Item(38,14).Value = 21000
$ws.Cells.Item(38,15).Value = 22000
$ws.Cells.Item(38,16).Value = 21500
$ws.Cells.Item(38,17).Value = "$/caja 20 kilos empedrada"
$ws.Cells.Item(38,18).Value = "Región de O'Higgins"
$ws.Cells.Item(38,19).Value = 1075
$ws.Cells.Item(38,20).Value = 20

# Row 39 : new "September Sun" record
$ws.Cells.Item(39,1).Value = 1
$ws.Cells.Item(39,2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(39,3).Value = "Arica y Parinacota"
$ws.Cells.Item(39,4).Value = 44615
$ws.Cells.Item(39,5).Value = 15
$ws.Cells.Item(39,6).Value = "Fruta"
$ws.Cells.Item(39,7).Value = 100103
$ws.Cells.Item(39,8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(39,9).Value = 100103004
$ws.Cells.Item(39,10).Value = "Durazno"
$ws.Cells.Item(39,11).Value = "September Sun"
$ws.Cells.Item(39,12).Value = "Segunda"
$ws.Cells.Item(39,13).Value = 300
$ws.Cells.Item(39,14).Value = 21000
$ws.Cells.Item(39,15).Value = 22000
$ws.Cells.Item(39,16).Value = 21500
$ws.Cells.Item(39,17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(39,18).Value = "Región de O'Higgins"
$ws.Cells.Item(39,19).Value = 1194
$ws.Cells.Item(39,20).Value = 18
